# Extend the "Days" log by 11 more days (rows 307-317), continuing the
# existing sequence-number / date-serial pattern, and update the active
# selection to match where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 306
$firstNewRow = $lastRow + 1
$lastNewRow = 317

# Copy the formatting of the last existing data row down onto the new
# rows (this also fills in the correct cell styles, borders, number
# format, etc. - matching style indices 3 (col A) / 4 (col B) used by
# every prior row) before writing the new values.
$ws.Range("A$lastRow`:B$lastRow").Copy()
$ws.Range("A$firstNewRow`:B$lastNewRow").PasteSpecial(-4122)

for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 2).Value = 43405 + ($r - $lastRow)
}

# Move the selection to where it ended up in the edited workbook.
$ws.Range("C315").Select()
